# "fix book type add book" - add four new book rows (10-13) to the Books sheet,
# reusing shared strings where the text already exists and introducing new
# shared strings ("aa", "akkaka", "rrea", "Romance", "koty") where it does not.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: The Great Fire / Fantasy / Super książka
$ws.Range("A10").Value = "The Great Fire"
$ws.Range("B10").Value = 37723.0
$ws.Range("D10").Value = "Fantasy"
$ws.Range("E10").Value = 345.0
$ws.Range("F10").Value = "Super książka"

# Row 11: Londyn NW / Horror / Całkiem niezła
$ws.Range("A11").Value = "Londyn NW"
$ws.Range("B11").Value = 40474.0
$ws.Range("D11").Value = "Horror"
$ws.Range("E11").Value = 343.0
$ws.Range("F11").Value = "Całkiem niezła"

# Row 12: aa / Horror / akkaka (new book, fixed type)
$ws.Range("A12").Value = "aa"
$ws.Range("B12").Value = -1.0
$ws.Range("D12").Value = "Horror"
$ws.Range("E12").Value = 112.0
$ws.Range("F12").Value = "akkaka"

# Row 13: rrea / Romance / koty (new book, fixed type)
$ws.Range("A13").Value = "rrea"
$ws.Range("B13").Value = -1.0
$ws.Range("D13").Value = "Romance"
$ws.Range("E13").Value = 12332990.0
$ws.Range("F13").Value = "koty"

# Columns D (type) and E (pages) need to grow to fit the new content, same as
# Excel's "bestFit" auto-resize that fires when text no longer fits the column.
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
